$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 127
$ws.Range("A127").Value = 111785216
$ws.Range("B127").Value = 78578
$ws.Range("D127").Value = 'NT'
$ws.Range("E127").Value = 6458
$ws.Range("F127").Value = 'Lunglav'
$ws.Range("G127").Value = 'Lobaria pulmonaria'
$ws.Range("H127").Value = '(L.) Hoffm.'
$ws.Range("Q127").Value = 577619.0624429594
$ws.Range("R127").Value = 6944938.738972194

# Row 128
$ws.Range("A128").Value = 111785209
$ws.Range("B128").Value = 78578
$ws.Range("D128").Value = 'NT'
$ws.Range("E128").Value = 6458
$ws.Range("F128").Value = 'Lunglav'
$ws.Range("G128").Value = 'Lobaria pulmonaria'
$ws.Range("H128").Value = '(L.) Hoffm.'
$ws.Range("Q128").Value = 577708.845606568
$ws.Range("R128").Value = 6945178.823951898

# Row 129
$ws.Range("A129").Value = 111785226
$ws.Range("B129").Value = 78578
$ws.Range("D129").Value = 'NT'
$ws.Range("E129").Value = 6458
$ws.Range("F129").Value = 'Lunglav'
$ws.Range("G129").Value = 'Lobaria pulmonaria'
$ws.Range("H129").Value = '(L.) Hoffm.'
$ws.Range("Q129").Value = 577680.815564253
$ws.Range("R129").Value = 6944566.386431335

# Row 130
$ws.Range("A130").Value = 111785246
$ws.Range("B130").Value = 73634
$ws.Range("D130").Value = 'LC'
$ws.Range("E130").Value = 6426
$ws.Range("F130").Value = 'Kattfotslav'
$ws.Range("G130").Value = 'Felipes leucopellaeus'
$ws.Range("H130").Value = '(Ach.) Frisch & G.Thor'
$ws.Range("Q130").Value = 577818.8902223237
$ws.Range("R130").Value = 6944791.983443609

# Row 131
$ws.Range("A131").Value = 111785208
$ws.Range("B131").Value = 78578
$ws.Range("D131").Value = 'NT'
$ws.Range("E131").Value = 6458
$ws.Range("F131").Value = 'Lunglav'
$ws.Range("G131").Value = 'Lobaria pulmonaria'
$ws.Range("H131").Value = '(L.) Hoffm.'
$ws.Range("Q131").Value = 577539.8790421919
$ws.Range("R131").Value = 6945032.627663832

# Row 132
$ws.Range("A132").Value = 111785225
$ws.Range("B132").Value = 78578
$ws.Range("D132").Value = 'NT'
$ws.Range("E132").Value = 6458
$ws.Range("F132").Value = 'Lunglav'
$ws.Range("G132").Value = 'Lobaria pulmonaria'
$ws.Range("H132").Value = '(L.) Hoffm.'
$ws.Range("Q132").Value = 577742.0418335226
$ws.Range("R132").Value = 6944530.994174051

# Row 133
$ws.Range("A133").Value = 111785210
$ws.Range("B133").Value = 78578
$ws.Range("D133").Value = 'NT'
$ws.Range("E133").Value = 6458
$ws.Range("F133").Value = 'Lunglav'
$ws.Range("G133").Value = 'Lobaria pulmonaria'
$ws.Range("H133").Value = '(L.) Hoffm.'
$ws.Range("Q133").Value = 577699.3262563417
$ws.Range("R133").Value = 6945152.363841761

# Row 134
$ws.Range("A134").Value = 111785219
$ws.Range("B134").Value = 78578
$ws.Range("D134").Value = 'NT'
$ws.Range("E134").Value = 6458
$ws.Range("F134").Value = 'Lunglav'
$ws.Range("G134").Value = 'Lobaria pulmonaria'
$ws.Range("H134").Value = '(L.) Hoffm.'
$ws.Range("Q134").Value = 577815.6496847487
$ws.Range("R134").Value = 6944812.162249871

# Row 135
$ws.Range("A135").Value = 111785224
$ws.Range("B135").Value = 78578
$ws.Range("D135").Value = 'NT'
$ws.Range("E135").Value = 6458
$ws.Range("F135").Value = 'Lunglav'
$ws.Range("G135").Value = 'Lobaria pulmonaria'
$ws.Range("H135").Value = '(L.) Hoffm.'
$ws.Range("Q135").Value = 577789.1149903627
$ws.Range("R135").Value = 6944587.805691725

# Row 136
$ws.Range("A136").Value = 111785195
$ws.Range("B136").Value = 78605
$ws.Range("D136").Value = 'LC'
$ws.Range("E136").Value = 6462
$ws.Range("F136").Value = 'Stuplav'
$ws.Range("G136").Value = 'Nephroma bellum'
$ws.Range("H136").Value = '(Spreng.) Tuck.'
$ws.Range("Q136").Value = 577817.9552042313
$ws.Range("R136").Value = 6944616.105901928

# Row 137
$ws.Range("A137").Value = 111785243
$ws.Range("B137").Value = 96348
$ws.Range("D137").Value = 'VU'
$ws.Range("E137").Value = 220787
$ws.Range("F137").Value = 'Knärot'
$ws.Range("G137").Value = 'Goodyera repens'
$ws.Range("H137").Value = '(L.) R. Br.'
$ws.Range("Q137").Value = 577759.6338390541
$ws.Range("R137").Value = 6944547.520700022

# Row 138
$ws.Range("A138").Value = 111785241
$ws.Range("B138").Value = 96348
$ws.Range("D138").Value = 'VU'
$ws.Range("E138").Value = 220787
$ws.Range("F138").Value = 'Knärot'
$ws.Range("G138").Value = 'Goodyera repens'
$ws.Range("H138").Value = '(L.) R. Br.'
$ws.Range("Q138").Value = 577823.0208926643
$ws.Range("R138").Value = 6944675.151490607

# Row 139
$ws.Range("A139").Value = 111785220
$ws.Range("B139").Value = 78578
$ws.Range("D139").Value = 'NT'
$ws.Range("E139").Value = 6458
$ws.Range("F139").Value = 'Lunglav'
$ws.Range("G139").Value = 'Lobaria pulmonaria'
$ws.Range("H139").Value = '(L.) Hoffm.'
$ws.Range("Q139").Value = 577929.7593229595
$ws.Range("R139").Value = 6944625.642482976

# Row 140
$ws.Range("A140").Value = 111785248
$ws.Range("B140").Value = 89423
$ws.Range("D140").Value = 'NT'
$ws.Range("E140").Value = 5432
$ws.Range("F140").Value = 'Granticka'
$ws.Range("G140").Value = 'Porodaedalea chrysoloma'
$ws.Range("H140").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q140").Value = 577835.1553672029
$ws.Range("R140").Value = 6944747.252110518

# Row 141
$ws.Range("A141").Value = 111785211
$ws.Range("B141").Value = 78578
$ws.Range("D141").Value = 'NT'
$ws.Range("E141").Value = 6458
$ws.Range("F141").Value = 'Lunglav'
$ws.Range("G141").Value = 'Lobaria pulmonaria'
$ws.Range("H141").Value = '(L.) Hoffm.'
$ws.Range("Q141").Value = 577681.525001083
$ws.Range("R141").Value = 6945125.248796649

# Row 142
$ws.Range("A142").Value = 111785227
$ws.Range("B142").Value = 78578
$ws.Range("D142").Value = 'NT'
$ws.Range("E142").Value = 6458
$ws.Range("F142").Value = 'Lunglav'
$ws.Range("G142").Value = 'Lobaria pulmonaria'
$ws.Range("H142").Value = '(L.) Hoffm.'
$ws.Range("Q142").Value = 577609.3004002962
$ws.Range("R142").Value = 6944686.704950654

# Row 143
$ws.Range("A143").Value = 111785249
$ws.Range("B143").Value = 77515
$ws.Range("D143").Value = 'NT'
$ws.Range("E143").Value = 6425
$ws.Range("F143").Value = 'Garnlav'
$ws.Range("G143").Value = 'Alectoria sarmentosa'
$ws.Range("H143").Value = '(Ach.) Ach.'
$ws.Range("Q143").Value = 577733.3514479286
$ws.Range("R143").Value = 6944900.913506362

# Row 144
$ws.Range("A144").Value = 111785212
$ws.Range("B144").Value = 78578
$ws.Range("D144").Value = 'NT'
$ws.Range("E144").Value = 6458
$ws.Range("F144").Value = 'Lunglav'
$ws.Range("G144").Value = 'Lobaria pulmonaria'
$ws.Range("H144").Value = '(L.) Hoffm.'
$ws.Range("Q144").Value = 577594.6139770675
$ws.Range("R144").Value = 6945057.386468799

# Row 145
$ws.Range("A145").Value = 111785242
$ws.Range("B145").Value = 96348
$ws.Range("D145").Value = 'VU'
$ws.Range("E145").Value = 220787
$ws.Range("F145").Value = 'Knärot'
$ws.Range("G145").Value = 'Goodyera repens'
$ws.Range("H145").Value = '(L.) R. Br.'
$ws.Range("Q145").Value = 577789.7963988667
$ws.Range("R145").Value = 6944558.818226521

# Row 146
$ws.Range("A146").Value = 111785236
$ws.Range("B146").Value = 96348
$ws.Range("D146").Value = 'VU'
$ws.Range("E146").Value = 220787
$ws.Range("F146").Value = 'Knärot'
$ws.Range("G146").Value = 'Goodyera repens'
$ws.Range("H146").Value = '(L.) R. Br.'
$ws.Range("Q146").Value = 577572.4076091016
$ws.Range("R146").Value = 6944824.864356839

# Row 147
$ws.Range("A147").Value = 111785223
$ws.Range("B147").Value = 78578
$ws.Range("D147").Value = 'NT'
$ws.Range("E147").Value = 6458
$ws.Range("F147").Value = 'Lunglav'
$ws.Range("G147").Value = 'Lobaria pulmonaria'
$ws.Range("H147").Value = '(L.) Hoffm.'
$ws.Range("Q147").Value = 577811.1671387866
$ws.Range("R147").Value = 6944591.08636965

# Row 148
$ws.Range("A148").Value = 111785215
$ws.Range("B148").Value = 78578
$ws.Range("D148").Value = 'NT'
$ws.Range("E148").Value = 6458
$ws.Range("F148").Value = 'Lunglav'
$ws.Range("G148").Value = 'Lobaria pulmonaria'
$ws.Range("H148").Value = '(L.) Hoffm.'
$ws.Range("Q148").Value = 577573.3984224057
$ws.Range("R148").Value = 6944998.890149554

# Row 149
$ws.Range("A149").Value = 111785237
$ws.Range("B149").Value = 96348
$ws.Range("D149").Value = 'VU'
$ws.Range("E149").Value = 220787
$ws.Range("F149").Value = 'Knärot'
$ws.Range("G149").Value = 'Goodyera repens'
$ws.Range("H149").Value = '(L.) R. Br.'
$ws.Range("Q149").Value = 577541.6974019273
$ws.Range("R149").Value = 6945053.384041801

# Row 150
$ws.Range("A150").Value = 111785247
$ws.Range("B150").Value = 89419
$ws.Range("D150").Value = 'NT'
$ws.Range("E150").Value = 1204
$ws.Range("F150").Value = 'Gränsticka'
$ws.Range("G150").Value = 'Phellopilus nigrolimitatus'
$ws.Range("H150").Value = '(Romell) Niemelä, T.Wagner & M.Fisch.'
$ws.Range("Q150").Value = 577741.5837880005
$ws.Range("R150").Value = 6944884.07477704

# Row 151
$ws.Range("A151").Value = 111785239
$ws.Range("B151").Value = 96348
$ws.Range("D151").Value = 'VU'
$ws.Range("E151").Value = 220787
$ws.Range("F151").Value = 'Knärot'
$ws.Range("G151").Value = 'Goodyera repens'
$ws.Range("H151").Value = '(L.) R. Br.'
$ws.Range("Q151").Value = 577769.2469415551
$ws.Range("R151").Value = 6944844.675943938

# Row 152
$ws.Range("A152").Value = 111785194
$ws.Range("B152").Value = 78605
$ws.Range("D152").Value = 'LC'
$ws.Range("E152").Value = 6462
$ws.Range("F152").Value = 'Stuplav'
$ws.Range("G152").Value = 'Nephroma bellum'
$ws.Range("H152").Value = '(Spreng.) Tuck.'
$ws.Range("Q152").Value = 577827.9269310302
$ws.Range("R152").Value = 6944682.172251224

# Row 153
$ws.Range("A153").Value = 111785213
$ws.Range("B153").Value = 78578
$ws.Range("D153").Value = 'NT'
$ws.Range("E153").Value = 6458
$ws.Range("F153").Value = 'Lunglav'
$ws.Range("G153").Value = 'Lobaria pulmonaria'
$ws.Range("H153").Value = '(L.) Hoffm.'
$ws.Range("Q153").Value = 577597.0534766318
$ws.Range("R153").Value = 6945012.333558927

# Row 154
$ws.Range("A154").Value = 111785218
$ws.Range("B154").Value = 78578
$ws.Range("D154").Value = 'NT'
$ws.Range("E154").Value = 6458
$ws.Range("F154").Value = 'Lunglav'
$ws.Range("G154").Value = 'Lobaria pulmonaria'
$ws.Range("H154").Value = '(L.) Hoffm.'
$ws.Range("Q154").Value = 577777.4147107385
$ws.Range("R154").Value = 6944830.597507096

# Row 155
$ws.Range("A155").Value = 111785196
$ws.Range("B155").Value = 78605
$ws.Range("D155").Value = 'LC'
$ws.Range("E155").Value = 6462
$ws.Range("F155").Value = 'Stuplav'
$ws.Range("G155").Value = 'Nephroma bellum'
$ws.Range("H155").Value = '(Spreng.) Tuck.'
$ws.Range("Q155").Value = 577807.1996835115
$ws.Range("R155").Value = 6944602.962691978

# Row 156
$ws.Range("A156").Value = 111785193
$ws.Range("B156").Value = 78605
$ws.Range("D156").Value = 'LC'
$ws.Range("E156").Value = 6462
$ws.Range("F156").Value = 'Stuplav'
$ws.Range("G156").Value = 'Nephroma bellum'
$ws.Range("H156").Value = '(Spreng.) Tuck.'
$ws.Range("Q156").Value = 577575.9394238007
$ws.Range("R156").Value = 6944851.186223409

# Row 157
$ws.Range("A157").Value = 111785222
$ws.Range("B157").Value = 78578
$ws.Range("D157").Value = 'NT'
$ws.Range("E157").Value = 6458
$ws.Range("F157").Value = 'Lunglav'
$ws.Range("G157").Value = 'Lobaria pulmonaria'
$ws.Range("H157").Value = '(L.) Hoffm.'
$ws.Range("Q157").Value = 577817.9552042313
$ws.Range("R157").Value = 6944616.105901928

# Row 158
$ws.Range("A158").Value = 111785240
$ws.Range("B158").Value = 96348
$ws.Range("D158").Value = 'VU'
$ws.Range("E158").Value = 220787
$ws.Range("F158").Value = 'Knärot'
$ws.Range("G158").Value = 'Goodyera repens'
$ws.Range("H158").Value = '(L.) R. Br.'
$ws.Range("Q158").Value = 577827.3226169772
$ws.Range("R158").Value = 6944747.067874849

# Row 159
$ws.Range("A159").Value = 111785207
$ws.Range("B159").Value = 78578
$ws.Range("D159").Value = 'NT'
$ws.Range("E159").Value = 6458
$ws.Range("F159").Value = 'Lunglav'
$ws.Range("G159").Value = 'Lobaria pulmonaria'
$ws.Range("H159").Value = '(L.) Hoffm.'
$ws.Range("Q159").Value = 577591.1417373432
$ws.Range("R159").Value = 6944930.258641767

# Row 160
$ws.Range("A160").Value = 111785221
$ws.Range("B160").Value = 78578
$ws.Range("D160").Value = 'NT'
$ws.Range("E160").Value = 6458
$ws.Range("F160").Value = 'Lunglav'
$ws.Range("G160").Value = 'Lobaria pulmonaria'
$ws.Range("H160").Value = '(L.) Hoffm.'
$ws.Range("Q160").Value = 577882.1692875527
$ws.Range("R160").Value = 6944590.915173424

# Row 161
$ws.Range("A161").Value = 111785238
$ws.Range("B161").Value = 96348
$ws.Range("D161").Value = 'VU'
$ws.Range("E161").Value = 220787
$ws.Range("F161").Value = 'Knärot'
$ws.Range("G161").Value = 'Goodyera repens'
$ws.Range("H161").Value = '(L.) R. Br.'
$ws.Range("Q161").Value = 577598.5769079959
$ws.Range("R161").Value = 6945065.304574955
